$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1842.5952
$ws.Range("I40").Value = 1803
$ws.Range("K40").Value = 1803
$ws.Range("M40").Value = -1628

$ws.Range("H127").Value = 877.2857
$ws.Range("I127").Value = 877.2857
$ws.Range("K127").Value = 2631.8571
$ws.Range("M127").Value = 2328.1429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 650.25
$ws.Range("I4").Value = 650.25
$ws.Range("K4").Value = 650.25
$ws.Range("M4").Value = -534.25

$ws.Range("H5").Value = 298.5
$ws.Range("I5").Value = 298.5
$ws.Range("K5").Value = 298.5
$ws.Range("M5").Value = -186.5

$ws.Range("H32").Value = 3409353.2
$ws.Range("I32").Value = 3591277.8
$ws.Range("K32").Value = 3591277.8
$ws.Range("M32").Value = -3590990.8

$ws.Range("H43").Value = 29874
$ws.Range("I43").Value = 15000
$ws.Range("J43").Value = 34832
$ws.Range("K43").Value = 15000
$ws.Range("L43").Value = 34832
$ws.Range("M43").Value = -14687
$ws.Range("N43").Value = -35458

$ws.Range("H61").Value = 4246.615
$ws.Range("I61").Value = 3677.0476
$ws.Range("J61").Value = 6638.8
$ws.Range("K61").Value = 3677.0476
$ws.Range("L61").Value = 6638.8
$ws.Range("M61").Value = -3465.0476
$ws.Range("N61").Value = -7062.8

$ws.Range("H63").Value = 7439.1904
$ws.Range("I63").Value = 1049.6666
$ws.Range("K63").Value = 1049.6666
$ws.Range("M63").Value = -363.6666

$ws.Range("H66").Value = 7439.1904
$ws.Range("I66").Value = 1049.6666
$ws.Range("K66").Value = 5248.333000000001
$ws.Range("M66").Value = -1816.333000000001

$ws.Range("H74").Value = 378648.38
$ws.Range("I74").Value = 532958.3
$ws.Range("K74").Value = 532958.3
$ws.Range("M74").Value = -532084.3

$ws.Range("H77").Value = 378648.38
$ws.Range("I77").Value = 532958.3
$ws.Range("K77").Value = 2664791.5
$ws.Range("M77").Value = -2660423.5

$ws.Range("H136").Value = 4246.615
$ws.Range("I136").Value = 3677.0476
$ws.Range("J136").Value = 6638.8
$ws.Range("K136").Value = 11031.1428
$ws.Range("L136").Value = 19916.4
$ws.Range("M136").Value = -8481.1428
$ws.Range("N136").Value = -25016.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 298.5
$ws.Range("I4").Value = 298.5
$ws.Range("K4").Value = 298.5
$ws.Range("M4").Value = -183.5

$ws.Range("H22").Value = 1606
$ws.Range("I22").Value = 1699.8572
$ws.Range("J22").Value = 949
$ws.Range("K22").Value = 1699.8572
$ws.Range("L22").Value = 949
$ws.Range("M22").Value = -1526.8572
$ws.Range("N22").Value = -1295

$ws.Range("H134").Value = 5560.8486
$ws.Range("I134").Value = 3823.0908
$ws.Range("J134").Value = 9036.362999999999
$ws.Range("K134").Value = 11469.2724
$ws.Range("L134").Value = 27109.089
$ws.Range("M134").Value = -8934.2724
$ws.Range("N134").Value = -32179.089

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2534.5
$ws.Range("I16").Value = 1163.8572
$ws.Range("J16").Value = 3600.5557
$ws.Range("K16").Value = 1163.8572
$ws.Range("L16").Value = 3600.5557
$ws.Range("M16").Value = -876.8571999999999
$ws.Range("N16").Value = -4174.5557

$ws.Range("H22").Value = 1061.2858
$ws.Range("I22").Value = 314.75
$ws.Range("J22").Value = 2056.6667
$ws.Range("K22").Value = 314.75
$ws.Range("L22").Value = 2056.6667
$ws.Range("M22").Value = 35.25
$ws.Range("N22").Value = -2756.6667

$ws.Range("H31").Value = 4137.7837
$ws.Range("I31").Value = 2014.2858
$ws.Range("J31").Value = 5430.3477
$ws.Range("K31").Value = 2014.2858
$ws.Range("L31").Value = 5430.3477
$ws.Range("M31").Value = -1719.2858
$ws.Range("N31").Value = -6020.3477

$ws.Range("H34").Value = 4137.7837
$ws.Range("I34").Value = 2014.2858
$ws.Range("J34").Value = 5430.3477
$ws.Range("K34").Value = 2014.2858
$ws.Range("L34").Value = 5430.3477
$ws.Range("M34").Value = -1812.2858
$ws.Range("N34").Value = -5834.3477

$ws.Range("H99").Value = 1214.6111
$ws.Range("I99").Value = 1271.5834
$ws.Range("J99").Value = 1100.6666
$ws.Range("K99").Value = 1271.5834
$ws.Range("L99").Value = 1100.6666
$ws.Range("M99").Value = 226.4166
$ws.Range("N99").Value = -4096.6666

$ws.Range("H113").Value = 2534.5
$ws.Range("I113").Value = 1163.8572
$ws.Range("J113").Value = 3600.5557
$ws.Range("K113").Value = 1163.8572
$ws.Range("L113").Value = 3600.5557
$ws.Range("M113").Value = 1006.1428
$ws.Range("N113").Value = -7940.5557

$ws.Range("H126").Value = 1214.6111
$ws.Range("I126").Value = 1271.5834
$ws.Range("J126").Value = 1100.6666
$ws.Range("K126").Value = 3814.7502
$ws.Range("L126").Value = 3301.9998
$ws.Range("M126").Value = -1344.7502
$ws.Range("N126").Value = -8241.9998

$ws.Range("H132").Value = 49109.223
$ws.Range("I132").Value = 3639.2222
$ws.Range("J132").Value = 140049.22
$ws.Range("K132").Value = 10917.6666
$ws.Range("L132").Value = 420147.66
$ws.Range("M132").Value = -8387.6666
$ws.Range("N132").Value = -425207.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 107010.7
$ws.Range("J68").Value = 9315.076999999999
$ws.Range("L68").Value = 27945.231
$ws.Range("N68").Value = -29567.231

$ws.Range("H71").Value = 107010.7
$ws.Range("J71").Value = 9315.076999999999
$ws.Range("L71").Value = 83835.693
$ws.Range("N71").Value = -91947.693

$ws.Range("H97").Value = 527.5454999999999
$ws.Range("I97").Value = 516.5
$ws.Range("J97").Value = 540.8
$ws.Range("K97").Value = 1549.5
$ws.Range("L97").Value = 1622.4
$ws.Range("M97").Value = -1053.5
$ws.Range("N97").Value = -2614.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 84119
$ws.Range("I2").Value = 922.6667
$ws.Range("J2").Value = 333708
$ws.Range("K2").Value = 922.6667
$ws.Range("L2").Value = 333708
$ws.Range("M2").Value = -809.6667
$ws.Range("N2").Value = -333934

$ws.Range("H41").Value = 6296
$ws.Range("J41").Value = 8918.5
$ws.Range("L41").Value = 8918.5
$ws.Range("N41").Value = -9628.5

$ws.Range("H58").Value = 21023
$ws.Range("J58").Value = 21023
$ws.Range("L58").Value = 21023
$ws.Range("N58").Value = -21577

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9237.056
$ws.Range("I132").Value = 6627.3
$ws.Range("J132").Value = 12499.25
$ws.Range("K132").Value = 19881.9
$ws.Range("L132").Value = 37497.75
$ws.Range("M132").Value = -17351.9
$ws.Range("N132").Value = -42557.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3092.3845
$ws.Range("I136").Value = 2182.1765
$ws.Range("J136").Value = 4811.6665
$ws.Range("K136").Value = 6546.529500000001
$ws.Range("L136").Value = 14434.9995
$ws.Range("M136").Value = -3996.529500000001
$ws.Range("N136").Value = -19534.9995
